# Add "fullRNASEQ" as the purpose for all sample rows (2-18), and select
# the next block of cells below the data (reflecting a user about to add
# another run's data / preparer info).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "purpose" column (E) for data rows 2 through 18 from "S.GISH"
# to the newly introduced value "fullRNASEQ".
$ws.Range("E2:E18").Value = "fullRNASEQ"

# Reflect the new selection left behind in the saved view.
$ws.Range("D19:F24").Select()
